$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in row 22 (was blank, now a "proxy for individual heat" data row) ---
$ws.Range("A22").Value2 = "remind"
$ws.Range("B22").Value2 = "SSP2-Base"
$ws.Range("C22").Value2 = "Wind"
$ws.Range("D22").Value2 = "DK"
$ws.Range("E22").Value2 = "consumption|individual heat|individual heat use"
$ws.Range("F22").Value2 = "PJ"
$ws.Range("G22").Value2 = 53.7
$ws.Range("H22").Value2 = 70.1

# E22:E24 pick up the "commented" number-format style (same one already used by E18)
$commentedStyle = $ws.Range("E22:E24")
$commentedStyle.NumberFormat = "0.0000"
$commentedStyle.Font.Name = "Arial"

# --- New cell comment on E22, duplicating the note already on E18 ---
$noteText = "Data in source doesn" + [char]0x2019 + "t add up, and unit is unclear. Must be wrong?"
$comment = $ws.Range("E22").AddComment($noteText)

# --- Update the saved selection / active cell ---
[void]$ws.Range("H23").Select()
